$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Torque")
Write-Host ("Formula=" + $ws.Range("C21").Formula)
Write-Host ("FormulaV=" + $ws.Range("C21").Value2)
